$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels from v1 to v2
$ws.Cells.Item(1, 6).Value  = "corn-n-soy-v2-exp"
$ws.Cells.Item(1, 8).Value  = "onlySoy-v2-exp"
$ws.Cells.Item(1, 9).Value  = "Diff: corn-n-soy-v2-exp"
$ws.Cells.Item(1, 11).Value = "Diff: onlySoy-v2-exp"
$ws.Cells.Item(1, 12).Value = "Error %: corn-n-soy-v2-exp"
$ws.Cells.Item(1, 14).Value = "Error %: onlySoy-v2-exp"

# Corrected per-row counts for corn-n-soy (col F) and onlySoy (col H),
# fixing a double-count bug in the soy detector.
$newCounts = @{
  2 = @{ F = 122; H = 82 }
  3 = @{ F = 93; H = 107 }
  4 = @{ F = 156; H = 128 }
  5 = @{ F = 185; H = 158 }
  6 = @{ F = 303; H = 302 }
  7 = @{ F = 206; H = 171 }
  8 = @{ F = 197; H = 286 }
  9 = @{ F = 290; H = 255 }
  10 = @{ F = 241; H = 237 }
  11 = @{ F = 231; H = 226 }
  12 = @{ F = 153; H = 145 }
  13 = @{ F = 118; H = 94 }
  14 = @{ F = 128; H = 207 }
  15 = @{ F = 516; H = 502 }
  16 = @{ F = 309; H = 193 }
  17 = @{ F = 134; H = 71 }
  18 = @{ F = 552; H = 498 }
  19 = @{ F = 344; H = 294 }
  20 = @{ F = 117; H = 113 }
  21 = @{ F = 368; H = 231 }
  22 = @{ F = 153; H = 196 }
  23 = @{ F = 259; H = 280 }
  24 = @{ F = 809; H = 642 }
  25 = @{ F = 466; H = 437 }
  26 = @{ F = 1918; H = 1880 }
  27 = @{ F = 86; H = 98 }
  28 = @{ F = 159; H = 179 }
  29 = @{ F = 559; H = 353 }
  30 = @{ F = 595; H = 583 }
  31 = @{ F = 326; H = 279 }
  32 = @{ F = 146; H = 159 }
  33 = @{ F = 366; H = 362 }
  34 = @{ F = 898; H = 876 }
  35 = @{ F = 1206; H = 1208 }
  36 = @{ F = 364; H = 365 }
  37 = @{ F = 514; H = 517 }
}

foreach ($row in $newCounts.Keys) {
    $newF = $newCounts[$row].F
    $newH = $newCounts[$row].H

    # Manual Count (E) and onlyCorn count (G) are untouched by this fix.
    $E = [double]$ws.Cells.Item($row, 5).Value2
    $G = [double]$ws.Cells.Item($row, 7).Value2

    $ws.Cells.Item($row, 6).Value = $newF   # corn-n-soy-v2-exp
    $ws.Cells.Item($row, 8).Value = $newH   # onlySoy-v2-exp

    $diffF = $newF - $E
    $diffG = $G - $E
    $diffH = $newH - $E

    $ws.Cells.Item($row, 9).Value  = $diffF                 # Diff: corn-n-soy-v2-exp
    $ws.Cells.Item($row, 10).Value = $diffG                 # Diff: onlyCorn-v1-exp6 (unchanged formula)
    $ws.Cells.Item($row, 11).Value = $diffH                 # Diff: onlySoy-v2-exp

    $ws.Cells.Item($row, 12).Value = ($diffF / $E) * 100    # Error %: corn-n-soy-v2-exp
    $ws.Cells.Item($row, 13).Value = ($diffG / $E) * 100    # Error %: onlyCorn-v1-exp6 (unchanged)
    $ws.Cells.Item($row, 14).Value = ($diffH / $E) * 100    # Error %: onlySoy-v2-exp
}

Write-Host "done"